$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header cells for the accrual currency / exchange rate columns
$ws.Range("L1").Value = "Accrual currency"
$ws.Range("M1").Value = "accrual exchange rate"

# Row 2 data
$ws.Range("L2").Value = "INRA"
$ws.Range("M2").Value = 1.223

# Row 3 data
$ws.Range("L3").Value = "A"
$ws.Range("M3").Value = 2.33

# Update the selection to match the authored change
$ws.Range("L1:M3").Select()
